$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the two label cells from A13/A14 up into H2/H3
$ws.Range("H2").Value = $ws.Range("A13").Value2
$ws.Range("H3").Value = $ws.Range("A14").Value2

# Clear the now-unused rows 13 and 14
$ws.Range("A13").ClearContents()
$ws.Range("A14").ClearContents()

# Update the selection shown when the sheet is reopened
$ws.Range("H2:H3").Select()
